$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.867.28"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").Value = "2.573.77"
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +1.98%  "

$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.366"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "3.032.50"
$ws.Range("E13").Value = "  -4.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "61.789.93"
$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "2.580.64"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("E21").Value = "  -4.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("E29").Value = "  -2.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  -1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "332.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.606"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.116.90"
$ws.Range("E46").Value = "  +0.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
